$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44635
$ws.Range("D3").Value = 44243
$ws.Range("I3").Value = "Especial"
$ws.Range("J3").Value = 300
$ws.Range("D4").Value = 44243
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = 10000
$ws.Range("P4").Value = 556
$ws.Range("D5").Value = 44243
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 8000
$ws.Range("L5").Value = 8000
$ws.Range("M5").Value = 8000
$ws.Range("P5").Value = 444
$ws.Range("D6").Value = 44238
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 12000
$ws.Range("L6").Value = 12000
$ws.Range("M6").Value = 12000
$ws.Range("P6").Value = 667
$ws.Range("D7").Value = 44238
$ws.Range("I7").Value = "Segunda"
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 10000
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = 10000
$ws.Range("P7").Value = 556
$ws.Range("D8").Value = 44238
$ws.Range("I8").Value = "Tercera"
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 8000
$ws.Range("L8").Value = 8000
$ws.Range("M8").Value = 8000
$ws.Range("P8").Value = 444
$ws.Range("D10").Value = 44585
$ws.Range("J10").Value = 200
$ws.Range("D11").Value = 44391
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 400
$ws.Range("K11").Value = 15000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 15000
$ws.Range("P11").Value = 833
$ws.Range("D12").Value = 44627
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 15000
$ws.Range("P12").Value = 833
$ws.Range("D13").Value = 44614
$ws.Range("N13").Value = "$/caja 18 kilos granel"
$ws.Range("D14").Value = 44628
$ws.Range("J14").Value = 300
$ws.Range("D15").Value = 44596
$ws.Range("J15").Value = 150
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 14000
$ws.Range("M15").Value = 14000
$ws.Range("P15").Value = 778
$ws.Range("D16").Value = 44245
$ws.Range("K16").Value = 12000
$ws.Range("L16").Value = 12000
$ws.Range("M16").Value = 12000
$ws.Range("P16").Value = 667
$ws.Range("D17").Value = 44245
$ws.Range("K17").Value = 10000
$ws.Range("L17").Value = 10000
$ws.Range("M17").Value = 10000
$ws.Range("P17").Value = 556
$ws.Range("D18").Value = 44396
$ws.Range("J18").Value = 250
$ws.Range("D19").Value = 44396
$ws.Range("I19").Value = "Segunda"
$ws.Range("J19").Value = 150
$ws.Range("D20").Value = 44235
$ws.Range("J20").Value = 400
$ws.Range("K20").Value = 13000
$ws.Range("L20").Value = 13000
$ws.Range("M20").Value = 13000
$ws.Range("P20").Value = 722
$ws.Range("D21").Value = 44235
$ws.Range("J21").Value = 200
$ws.Range("K21").Value = 11000
$ws.Range("L21").Value = 11000
$ws.Range("M21").Value = 11000
$ws.Range("P21").Value = 611
$ws.Range("D22").Value = 44235
$ws.Range("I22").Value = "Tercera"
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 9000
$ws.Range("L22").Value = 9000
$ws.Range("M22").Value = 9000
$ws.Range("P22").Value = 500
$ws.Range("D23").Value = 44630
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 15000
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = 15000
$ws.Range("P23").Value = 833
$ws.Range("D26").Value = 44383
$ws.Range("K26").Value = 16000
$ws.Range("L26").Value = 16000
$ws.Range("M26").Value = 16000
$ws.Range("N26").Value = "$/bandeja 18 kilos"
$ws.Range("P26").Value = 889
$ws.Range("D27").Value = 44383
$ws.Range("I27").Value = "Segunda"
$ws.Range("J27").Value = 200
$ws.Range("D28").Value = 44631
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 300
$ws.Range("K28").Value = 15000
$ws.Range("L28").Value = 15000
$ws.Range("M28").Value = 15000
$ws.Range("P28").Value = 833